# Edit script for "Algoritmo de Busqueda en Profundidad (DFS).docx"
#
# 1) Splits the title run "... (DFS: Depth First Search)" into separate
#    runs so that the English words "Depth", "First" and "Search" are each
#    wrapped in a spell-check proofErr (spellStart/spellEnd) pair, matching
#    what Word's background spell checker produces for foreign words.
# 2) Turns on underline for the paragraph mark of the
#    "Cantidad de Ramas del Recorrido" paragraph and appends
#    " - Nodo Inicial (Raiz del Arbol del Recorrido)" as extra text.
#
# Both changes are applied by locating the target Range and replacing its
# contents with equivalent WordprocessingML via Range.InsertXML, which lets
# us control run splitting and proofErr markers precisely.

$d = $word.ActiveDocument

# --- Change 1: title paragraph -------------------------------------------------
$titleText = "Algoritmo de Búsqueda en Profundidad (DFS: Depth First Search)"

$titleRange = $d.Content
$found = $titleRange.Find.Execute($titleText, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the title text to update"
}

$titleXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:color w:val="548DD4" w:themeColor="text2" w:themeTint="99"/><w:sz w:val="44"/></w:rPr><w:t xml:space="preserve">Algoritmo de Búsqueda en Profundidad (DFS: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="548DD4" w:themeColor="text2" w:themeTint="99"/><w:sz w:val="44"/></w:rPr><w:t>Depth</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="548DD4" w:themeColor="text2" w:themeTint="99"/><w:sz w:val="44"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="548DD4" w:themeColor="text2" w:themeTint="99"/><w:sz w:val="44"/></w:rPr><w:t>First</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="548DD4" w:themeColor="text2" w:themeTint="99"/><w:sz w:val="44"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="548DD4" w:themeColor="text2" w:themeTint="99"/><w:sz w:val="44"/></w:rPr><w:t>Search</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="548DD4" w:themeColor="text2" w:themeTint="99"/><w:sz w:val="44"/></w:rPr><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titleRange.InsertXML($titleXml)

# --- Change 2: "Cantidad de Ramas del Recorrido" paragraph ---------------------
$branchPara = $null
foreach ($p in $d.Paragraphs) {
    $paraText = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($paraText -eq "Cantidad de Ramas del Recorrido") {
        $branchPara = $p
        break
    }
}
if ($null -eq $branchPara) {
    throw "Could not find the 'Cantidad de Ramas del Recorrido' paragraph"
}

$branchXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00FF130C" w:rsidRDefault="00012A65"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>Cantidad de Ramas del Recorrido</w:t></w:r><w:r><w:t xml:space="preserve"> – Nodo Inicial (Raíz del Árbol del Recorrido)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$branchPara.Range.InsertXML($branchXml)
